# Update workbook to suit the new layout of MS Template.

$wb = $excel.ActiveWorkbook

# --- Rename sheets to the new MS Template layout ---
$wb.Worksheets.Item("ISTD map").Name          = "Transition_Name_Annot"
$wb.Worksheets.Item("normArea by ISTD").Name  = "normArea_by_ISTD"
$wb.Worksheets.Item("Sample Annot").Name      = "Sample_Annot"
$wb.Worksheets.Item("normConc by ISTD").Name  = "normConc_by_ISTD"

# --- Update the ISTD concentration unit header (nM -> uM) ---
$istdSheet = $wb.Worksheets.Item("Transition_Name_Annot")
$istdUsed = $istdSheet.UsedRange
for ($c = 1; $c -le $istdUsed.Columns.Count; $c++) {
    $cell = $istdSheet.Cells.Item(1, $c)
    if ($cell.Value2 -eq "ISTD_Conc_[nM]") {
        $cell.Value = "ISTD_Conc_[uM]"
    }
}

# --- Rename Sample_Type category labels: PQC -> BQC, Sample -> SPL ---
$annotSheet = $wb.Worksheets.Item("Sample_Annot")
$annotUsed = $annotSheet.UsedRange
$firstRow = $annotUsed.Row
$lastRow = $firstRow + $annotUsed.Rows.Count - 1
$firstCol = $annotUsed.Column
$lastCol = $firstCol + $annotUsed.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $annotSheet.Cells.Item($r, $c)
        if ($cell.Value2 -eq "PQC") {
            $cell.Value = "BQC"
        } elseif ($cell.Value2 -eq "Sample") {
            $cell.Value = "SPL"
        }
    }
}
